$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.250998973846436
$ws.Range("B1").Value = 3.208706140518188
$ws.Range("C1").Value = 6.063765048980713
$ws.Range("D1").Value = 1.795663952827454
$ws.Range("E1").Value = 1.05532443523407
